# Update the "Förändrad" date column (C) for rows 2-28: increment each
# date serial value by 1 day (45534 -> 45535), matching the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = $cell.Value2 + 1
}
